$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 360; this shifts the existing rows 360-373
# (three weeks of "Comercializadora del Agro de Limari - Limon" data) down
# to rows 363-376, preserving all of their values/format.
$ws.Rows("360:362").Insert()

# Populate the 3 newly inserted rows with a new week's data (Fecha 44448),
# following the same boilerplate columns as the rest of the block.
$newRows = @(
    @{ Row=360; L="1a amarillo"; M=800; N=2800; O=3000; P=2900; S=181 },
    @{ Row=361; L="2a amarillo"; M=620; N=1800; O=2000; P=1900; S=119 },
    @{ Row=362; L="3a amarillo"; M=470; N=800;  O=1000; P=900;  S=56  }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 2
    $ws.Cells.Item($row, 2).Value = "Comercializadora del Agro de Limarí"
    $ws.Cells.Item($row, 3).Value = "Coquimbo"
    $ws.Cells.Item($row, 4).Value = 44448
    $ws.Cells.Item($row, 5).Value = 4
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100102
    $ws.Cells.Item($row, 8).Value = "Cítricos"
    $ws.Cells.Item($row, 9).Value = 100102003
    $ws.Cells.Item($row, 10).Value = "Limón"
    $ws.Cells.Item($row, 11).Value = "Sin especificar"
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = '$/malla 16 kilos'
    $ws.Cells.Item($row, 18).Value = "Provincia de Limarí"
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = 16
}
